$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.929813666666668
$ws.Range("H2").Value = 26.789441
$ws.Range("I2").Value = 0.3579859341865942
$ws.Range("J2").Value = 0.3579859341865942
$ws.Range("M2").Value = 4.993165333333334
$ws.Range("N2").Value = 14.979496
$ws.Range("O2").Value = 0.06779298131037136
$ws.Range("P2").Value = 0.06779298131037137
$ws.Range("Q2").Value = 44.58803603352623
$ws.Range("R2").Value = 401.2923243017361
$ws.Range("S2").Value = 0.02426893374568761
$ws.Range("T2").Value = 0.02426893374568762
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.929813666666668
$ws.Range("H3").Value = 26.789441
$ws.Range("I3").Value = 0.3579859341865942
$ws.Range("J3").Value = 0.3579859341865942
$ws.Range("O3").Value = 0.5355771637189464
$ws.Range("P3").Value = 0.5355771637189464
$ws.Range("Q3").Value = 352.2537792711117
$ws.Range("R3").Value = 3170.284013440005
$ws.Range("S3").Value = 0.1917290912829335
$ws.Range("T3").Value = 0.1917290912829335
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.929813666666668
$ws.Range("H4").Value = 26.789441
$ws.Range("I4").Value = 0.3579859341865942
$ws.Range("J4").Value = 0.3579859341865942
$ws.Range("M4").Value = 29.08216166666666
$ws.Range("N4").Value = 87.24648499999999
$ws.Range("O4").Value = 0.3948530262300277
$ws.Range("P4").Value = 0.3948530262300277
$ws.Range("Q4").Value = 259.6982847072094
$ws.Range("R4").Value = 2337.284562364885
$ws.Range("S4").Value = 0.1413518294613602
$ws.Range("T4").Value = 0.1413518294613602
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.929813666666668
$ws.Range("H5").Value = 26.789441
$ws.Range("I5").Value = 0.3579859341865942
$ws.Range("J5").Value = 0.3579859341865942
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.130869
$ws.Range("N5").Value = 0.392607
$ws.Range("O5").Value = 0.001776828740654623
$ws.Range("P5").Value = 0.001776828740654624
$ws.Range("Q5").Value = 1.168635784743
$ws.Range("R5").Value = 10.517722062687
$ws.Range("S5").Value = 0.000636079696612835
$ws.Range("T5").Value = 0.0006360796966128351
$ws.Range("I6").Value = 0.03080543264277933
$ws.Range("J6").Value = 0.03080543264277933
$ws.Range("M6").Value = 4.993165333333334
$ws.Range("N6").Value = 14.979496
$ws.Range("O6").Value = 0.06779298131037136
$ws.Range("P6").Value = 0.06779298131037137
$ws.Range("Q6").Value = 3.836893043928
$ws.Range("R6").Value = 34.532037395352
$ws.Range("S6").Value = 0.002088392119409843
$ws.Range("T6").Value = 0.002088392119409844
$ws.Range("I7").Value = 0.03080543264277933
$ws.Range("J7").Value = 0.03080543264277933
$ws.Range("O7").Value = 0.5355771637189464
$ws.Range("P7").Value = 0.5355771637189464
$ws.Range("S7").Value = 0.0164986862419548
$ws.Range("T7").Value = 0.0164986862419548
$ws.Range("I8").Value = 0.03080543264277933
$ws.Range("J8").Value = 0.03080543264277933
$ws.Range("M8").Value = 29.08216166666666
$ws.Range("N8").Value = 87.24648499999999
$ws.Range("O8").Value = 0.3948530262300277
$ws.Range("P8").Value = 0.3948530262300277
$ws.Range("Q8").Value = 22.347576407355
$ws.Range("R8").Value = 201.128187666195
$ws.Range("S8").Value = 0.0121636183033267
$ws.Range("T8").Value = 0.0121636183033267
$ws.Range("I9").Value = 0.03080543264277933
$ws.Range("J9").Value = 0.03080543264277933
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.130869
$ws.Range("N9").Value = 0.392607
$ws.Range("O9").Value = 0.001776828740654623
$ws.Range("P9").Value = 0.001776828740654624
$ws.Range("Q9").Value = 0.100563534801
$ws.Range("R9").Value = 0.905071813209
$ws.Range("S9").Value = 0.00005473597808799043
$ws.Range("T9").Value = 0.00005473597808799044
$ws.Range("G10").Value = 15.246351
$ws.Range("H10").Value = 45.739053
$ws.Range("I10").Value = 0.6112086331706265
$ws.Range("J10").Value = 0.6112086331706265
$ws.Range("M10").Value = 4.993165333333334
$ws.Range("N10").Value = 14.979496
$ws.Range("O10").Value = 0.06779298131037136
$ws.Range("P10").Value = 0.06779298131037137
$ws.Range("Q10").Value = 76.127551273032
$ws.Range("R10").Value = 685.147961457288
$ws.Range("S10").Value = 0.0414356554452739
$ws.Range("T10").Value = 0.04143565544527391
$ws.Range("G11").Value = 15.246351
$ws.Range("H11").Value = 45.739053
$ws.Range("I11").Value = 0.6112086331706265
$ws.Range("J11").Value = 0.6112086331706265
$ws.Range("O11").Value = 0.5355771637189464
$ws.Range("P11").Value = 0.5355771637189464
$ws.Range("Q11").Value = 601.421816884185
$ws.Range("R11").Value = 5412.796351957665
$ws.Range("S11").Value = 0.3273493861940581
$ws.Range("T11").Value = 0.3273493861940581
$ws.Range("G12").Value = 15.246351
$ws.Range("H12").Value = 45.739053
$ws.Range("I12").Value = 0.6112086331706265
$ws.Range("J12").Value = 0.6112086331706265
$ws.Range("M12").Value = 29.08216166666666
$ws.Range("N12").Value = 87.24648499999999
$ws.Range("O12").Value = 0.3948530262300277
$ws.Range("P12").Value = 0.3948530262300277
$ws.Range("Q12").Value = 443.3968446087449
$ws.Range("R12").Value = 3990.571601478704
$ws.Range("S12").Value = 0.2413375784653407
$ws.Range("T12").Value = 0.2413375784653408
$ws.Range("G13").Value = 15.246351
$ws.Range("H13").Value = 45.739053
$ws.Range("I13").Value = 0.6112086331706265
$ws.Range("J13").Value = 0.6112086331706265
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.130869
$ws.Range("N13").Value = 0.392607
$ws.Range("O13").Value = 0.001776828740654623
$ws.Range("P13").Value = 0.001776828740654624
$ws.Range("Q13").Value = 1.995274709019
$ws.Range("R13").Value = 17.957472381171
$ws.Range("S13").Value = 0.001086013065953798
$ws.Range("T13").Value = 0.001086013065953798
